# "before delete the comparison with LLA"
#
# The underlying data table on Sheet1 (A2:A13 - the x-axis "Deadline
# (cycles)" values feeding the stacked bar chart) is shifted by +2 for
# every row, e.g. 28,30,32,...,50 becomes 30,32,34,...,52. The chart's
# cached category values mirror the worksheet cells and move in lock-step.
# The sheet's active selection also moved from T19 to D17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Bump each value in A2:A13 by 2. Use Value2 (not Value) when reading back
# the current contents - Value round-trips as a descriptor string in this
# host, which would coerce the cell into text instead of a number.
for ($r = 2; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 2
}

# Move the active selection to D17 (was T19).
$ws.Range("D17").Select()
